$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.952.17"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "3.487.75"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.98"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.46"
$ws.Range("E6").Value = "  +3.19%  "

$ws.Range("E7").Value = "  +5.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "3.485.45"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("E10").Value = "  +8.00%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Value = "4.083.92"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.16"
$ws.Range("E14").Value = "  +1.60%  "

$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "67.890.65"
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").Value = "3.481.87"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.17"
$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.90"
$ws.Range("E21").Value = "  +1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.30"
$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.39"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("E33").Value = "  +0.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.68"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.36"
$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.72"
$ws.Range("E38").Value = "  -1.17%  "

$ws.Range("E39").Value = "  +3.14%  "

$ws.Range("E40").Value = "  +12.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.73"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.31"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0720"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "2.754.25"
$ws.Range("E46").Value = "  -1.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.36"
$ws.Range("E47").Value = "  -3.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.65"
$ws.Range("E48").Value = "  -1.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0300"
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "329.65"
$ws.Range("E50").Value = "  -3.55%  "

$ws.Range("E51").Value = "  -2.08%  "
